# MasterQuest.xlsx update:
#  - Updating the dialogue response quest + localization for quest descriptions
#  - Rewrites the "comment" column (Y) on the Entities sheet:
#      * Y1 note gets a trailing period added
#      * Y2 gets a brand-new note about filling in taskRequiredAmount
#      * Y3's old note ("Green is main quest...") is removed (cell left blank)
#      * The note that used to live in Y2 ("rewardKey that is empty...")
#        is moved down to the newly added Y4 row
#      * Y2:Y4 get a new "Note w/ colored font" style
#  - Updates the window selection on the Entities sheet

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Entities")

# ----- Update / move the note texts in column Y -----
$ws.Range("Y1").Value2 = "taskActionType of DialogueResponse need to fill in dialogue id into taskObjectiveKey."
$ws.Range("Y2").Value2 = "fill in the answer in taskRequiredAmount."
$ws.Range("Y4").Value2 = "rewardKey that is empty adds toward story progression."

# ----- Give Y2 the "Note" cell style with the reddish Input font color -----
$ws.Range("Y2").Style = "Note"
$ws.Range("Y2").Font.Color = 7749439

# ----- Copy that exact formatting onto Y3 and Y4 (reuses the same style, -----
# ----- avoids generating duplicate/unused style entries) -----
$ws.Range("Y2").Copy() | Out-Null
$ws.Range("Y3").PasteSpecial(-4122) | Out-Null
$ws.Range("Y4").PasteSpecial(-4122) | Out-Null

# Y3 no longer holds any text - the "Green is main quest..." note was removed
$ws.Range("Y3").Value2 = $null

# ----- Update the active window view/selection -----
$ws.Activate() | Out-Null
$win = $excel.ActiveWindow
$win.ScrollColumn = 11
$ws.Range("Y11").Select() | Out-Null
